$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.126374840736389
$ws.Range("B1").Value = 2.275492429733276
$ws.Range("C1").Value = 10.4202241897583
$ws.Range("D1").Value = 1.98122501373291
$ws.Range("E1").Value = 1.28537118434906
